$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.681.41'
$ws.Range("E2").Value = '  -5.54%  '
$ws.Range("D3").Value = '2.617.04'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''302.09'
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("D6").Value = '''96.18'
$ws.Range("E6").Value = '  -4.01%  '
$ws.Range("D7").Value = '''0.581'
$ws.Range("E7").Value = '  -3.79%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '''0.556'
$ws.Range("E9").Value = '  -4.80%  '
$ws.Range("D10").Value = '''36.84'
$ws.Range("E10").Value = '  -6.84%  '
$ws.Range("D11").Value = '''0.0813'
$ws.Range("E11").Value = '  -3.87%  '
$ws.Range("D12").Value = '''7.81'
$ws.Range("E12").Value = '  -5.06%  '
$ws.Range("D13").Value = '3.030.51'
$ws.Range("E13").Value = '  +1.11%  '
$ws.Range("E14").Value = '  +1.09%  '
$ws.Range("D15").Value = '2.622.54'
$ws.Range("E15").Value = '  -1.42%  '
$ws.Range("D16").Value = '''0.889'
$ws.Range("E16").Value = '  -3.85%  '
$ws.Range("D17").Value = '''14.38'
$ws.Range("E17").Value = '  -4.34%  '
$ws.Range("D18").Value = '43.660.98'
$ws.Range("E18").Value = '  -5.95%  '
$ws.Range("D19").Value = '''6.65'
$ws.Range("E19").Value = '  -1.50%  '
$ws.Range("D20").Value = '0.0₃0976'
$ws.Range("E20").Value = '  -3.80%  '
$ws.Range("D21").Value = '''12.50'
$ws.Range("E21").Value = '  -4.46%  '
$ws.Range("D22").Value = '''73.37'
$ws.Range("E22").Value = '  +2.15%  '
$ws.Range("D23").Value = '''267.33'
$ws.Range("E23").Value = '  -2.21%  '
$ws.Range("D24").Value = '''2.94'
$ws.Range("E24").Value = '  -3.17%  '
$ws.Range("D25").Value = '''2.22'
$ws.Range("E25").Value = '  +2.09%  '
$ws.Range("D26").Value = '''29.42'
$ws.Range("E26").Value = '  -2.73%  '
$ws.Range("D27").Value = '''1.00'
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").Value = '''10.24'
$ws.Range("E28").Value = '  -3.63%  '
$ws.Range("E29").Value = '  -3.76%  '
$ws.Range("D30").Value = '''37.96'
$ws.Range("E30").Value = '  -3.87%  '
$ws.Range("E31").Value = '  -3.75%  '
$ws.Range("D32").Value = '''3.62'
$ws.Range("E32").Value = '  -0.96%  '
$ws.Range("D33").Value = '''2.24'
$ws.Range("E33").Value = '  +0.43%  '
$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").Value = '''152.38'
$ws.Range("E34").Value = '  +1.44%  '
$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '''2.80'
$ws.Range("E35").Value = '  -1.46%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '''0.0810'
$ws.Range("E36").Value = '  -3.95%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '''0.117'
$ws.Range("E37").Value = '  -4.13%  '
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").Value = '''25.02'
$ws.Range("E38").Value = '  +7.01%  '
$ws.Range("E39").Value = '  -1.94%  '
$ws.Range("D40").Value = '''16.32'
$ws.Range("E40").Value = '  +1.15%  '
$ws.Range("D41").Value = '''3.49'
$ws.Range("E41").Value = '  -4.14%  '
$ws.Range("D42").Value = '''0.0316'
$ws.Range("E42").Value = '  -4.58%  '
$ws.Range("D43").Value = '''3.85'
$ws.Range("E43").Value = '  -6.49%  '
$ws.Range("D44").Value = '2.092.63'
$ws.Range("E44").Value = '  -3.35%  '
$ws.Range("D45").Value = '''0.997'
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").Value = '''88.69'
$ws.Range("E46").Value = '  -5.58%  '
$ws.Range("D47").Value = '''9.17'
$ws.Range("E47").Value = '  -4.89%  '
$ws.Range("D48").Value = '2.881.38'
$ws.Range("E48").Value = '  +0.92%  '
$ws.Range("D49").Value = '''1.61'
$ws.Range("E49").Value = '  +4.31%  '
$ws.Range("D50").Value = '''106.56'
$ws.Range("E50").Value = '  -2.82%  '
$ws.Range("D51").Value = '''0.191'
$ws.Range("E51").Value = '  -4.72%  '
